$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("z")

$data = @(
    @("EU27", "Activity", "Production of electricity by biomass and waste", "EU27", "Commodity", "Electricity", "Update", 0),
    @("EU27", "Activity", "Production of electricity by coal", "EU27", "Commodity", "Electricity", "Update", 0),
    @("EU27", "Activity", "Production of electricity by gas", "EU27", "Commodity", "Electricity", "Update", 0.02377972465581978),
    @("EU27", "Activity", "Production of electricity by hydro", "EU27", "Commodity", "Electricity", "Update", 0.04755944931163955),
    @("EU27", "Activity", "Production of electricity by nuclear", "EU27", "Commodity", "Electricity", "Update", 0),
    @("EU27", "Activity", "Production of electricity by petroleum and other oil derivatives", "EU27", "Commodity", "Electricity", "Update", 0),
    @("EU27", "Activity", "Production of electricity by solar photovoltaic", "EU27", "Commodity", "Electricity", "Update", 0.1902377972465582),
    @("EU27", "Activity", "Production of electricity by wind", "EU27", "Commodity", "Electricity", "Update", 0.7384230287859824)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $row++
}
